$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet (tab name "Sheet" -> "Sheet1")
$ws.Name = "Sheet1"

# Update row 2 content
$ws.Range("A2").Value = "Item1"
$ws.Range("B2").Value = 100

# Add new row 3 content
$ws.Range("A3").Value = "Item2"
$ws.Range("B3").Value = 200
